$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(428, 1).Value = "AY352487"
$ws.Cells.Item(428, 2).Value = 0.74099999999999999
$ws.Cells.Item(428, 3).Value = 1.2589999999999999
$ws.Cells.Item(428, 4).Value = 0.68600000000000005
$ws.Cells.Item(428, 5).Value = 1.5429999999999999
$ws.Cells.Item(428, 6).Value = 0.51400000000000001
$ws.Cells.Item(428, 7).Value = 0.51400000000000001
$ws.Cells.Item(428, 8).Value = 0.85699999999999998
$ws.Cells.Item(428, 9).Value = 1.8859999999999999
$ws.Cells.Item(428, 10).Value = 1
$ws.Cells.Item(428, 11).Value = 1
$ws.Cells.Item(428, 12).Value = 1
$ws.Cells.Item(428, 13).Value = 0.88900000000000001
$ws.Cells.Item(428, 14).Value = 1.333
$ws.Cells.Item(428, 15).Value = 0.88900000000000001
$ws.Cells.Item(428, 16).Value = 0.88900000000000001
$ws.Cells.Item(428, 17).Value = 1.2
$ws.Cells.Item(428, 18).Value = 0.85699999999999998
$ws.Cells.Item(428, 19).Value = 1.714
$ws.Cells.Item(428, 20).Value = 0.68600000000000005
$ws.Cells.Item(428, 21).Value = 0.68600000000000005
$ws.Cells.Item(428, 22).Value = 0.85699999999999998
$ws.Cells.Item(428, 23).Value = 1.647
$ws.Cells.Item(428, 24).Value = 1.1759999999999999
$ws.Cells.Item(428, 25).Value = 1.1759999999999999
$ws.Cells.Item(428, 26).Value = 0
$ws.Cells.Item(428, 27).Value = 1.103
$ws.Cells.Item(428, 28).Value = 1.103
$ws.Cells.Item(428, 29).Value = 0.55200000000000005
$ws.Cells.Item(428, 30).Value = 1.2410000000000001
$ws.Cells.Item(428, 31).Value = 1.4119999999999999
$ws.Cells.Item(428, 32).Value = 0.70599999999999996
$ws.Cells.Item(428, 33).Value = 1.7649999999999999
$ws.Cells.Item(428, 34).Value = 0.11799999999999999
$ws.Cells.Item(428, 35).Value = 1.429
$ws.Cells.Item(428, 36).Value = 0.57099999999999995
$ws.Cells.Item(428, 37).Value = 0.92300000000000004
$ws.Cells.Item(428, 38).Value = 1.077
$ws.Cells.Item(428, 39).Value = 1.2
$ws.Cells.Item(428, 40).Value = 0.8
$ws.Cells.Item(428, 41).Value = 1
$ws.Cells.Item(428, 42).Value = 1
$ws.Cells.Item(428, 43).Value = 0.53800000000000003
$ws.Cells.Item(428, 44).Value = 1.462
$ws.Cells.Item(428, 45).Value = 1.0429999999999999
$ws.Cells.Item(428, 46).Value = 0.95699999999999996
$ws.Cells.Item(428, 47).Value = 0.56200000000000006
$ws.Cells.Item(428, 48).Value = 1.4379999999999999
$ws.Cells.Item(428, 49).Value = 1.333
$ws.Cells.Item(428, 50).Value = 0.66700000000000004
$ws.Cells.Item(428, 51).Value = 0.78300000000000003
$ws.Cells.Item(428, 52).Value = 0
$ws.Cells.Item(428, 53).Value = 0.26100000000000001
$ws.Cells.Item(428, 54).Value = 0.78300000000000003
$ws.Cells.Item(428, 55).Value = 3.13
$ws.Cells.Item(428, 56).Value = 1.0429999999999999
$ws.Cells.Item(428, 57).Value = 0.64500000000000002
$ws.Cells.Item(428, 58).Value = 0.51600000000000001
$ws.Cells.Item(428, 59).Value = 1.032
$ws.Cells.Item(428, 60).Value = 1.806

$ws.Cells.Item(429, 1).Value = "EF611828"
$ws.Cells.Item(429, 2).Value = 0.74099999999999999
$ws.Cells.Item(429, 3).Value = 1.2589999999999999
$ws.Cells.Item(429, 4).Value = 0.68600000000000005
$ws.Cells.Item(429, 5).Value = 1.714
$ws.Cells.Item(429, 6).Value = 0.51400000000000001
$ws.Cells.Item(429, 7).Value = 0.51400000000000001
$ws.Cells.Item(429, 8).Value = 0.85699999999999998
$ws.Cells.Item(429, 9).Value = 1.714
$ws.Cells.Item(429, 10).Value = 1.111
$ws.Cells.Item(429, 11).Value = 0.88900000000000001
$ws.Cells.Item(429, 12).Value = 1
$ws.Cells.Item(429, 13).Value = 0.88900000000000001
$ws.Cells.Item(429, 14).Value = 1.333
$ws.Cells.Item(429, 15).Value = 0.88900000000000001
$ws.Cells.Item(429, 16).Value = 0.88900000000000001
$ws.Cells.Item(429, 17).Value = 1.0289999999999999
$ws.Cells.Item(429, 18).Value = 0.85699999999999998
$ws.Cells.Item(429, 19).Value = 1.714
$ws.Cells.Item(429, 20).Value = 0.85699999999999998
$ws.Cells.Item(429, 21).Value = 0.85699999999999998
$ws.Cells.Item(429, 22).Value = 0.68600000000000005
$ws.Cells.Item(429, 23).Value = 1.647
$ws.Cells.Item(429, 24).Value = 1.1759999999999999
$ws.Cells.Item(429, 25).Value = 1.1759999999999999
$ws.Cells.Item(429, 26).Value = 0
$ws.Cells.Item(429, 27).Value = 1.2
$ws.Cells.Item(429, 28).Value = 0.93300000000000005
$ws.Cells.Item(429, 29).Value = 0.66700000000000004
$ws.Cells.Item(429, 30).Value = 1.2
$ws.Cells.Item(429, 31).Value = 1.333
$ws.Cells.Item(429, 32).Value = 0.84799999999999998
$ws.Cells.Item(429, 33).Value = 1.6970000000000001
$ws.Cells.Item(429, 34).Value = 0.121
$ws.Cells.Item(429, 35).Value = 1.429
$ws.Cells.Item(429, 36).Value = 0.57099999999999995
$ws.Cells.Item(429, 37).Value = 0.92300000000000004
$ws.Cells.Item(429, 38).Value = 1.077
$ws.Cells.Item(429, 39).Value = 1.2
$ws.Cells.Item(429, 40).Value = 0.8
$ws.Cells.Item(429, 41).Value = 0.9
$ws.Cells.Item(429, 42).Value = 1.1000000000000001
$ws.Cells.Item(429, 43).Value = 0.53800000000000003
$ws.Cells.Item(429, 44).Value = 1.462
$ws.Cells.Item(429, 45).Value = 1.0429999999999999
$ws.Cells.Item(429, 46).Value = 0.95699999999999996
$ws.Cells.Item(429, 47).Value = 0.56200000000000006
$ws.Cells.Item(429, 48).Value = 1.4379999999999999
$ws.Cells.Item(429, 49).Value = 1.333
$ws.Cells.Item(429, 50).Value = 0.66700000000000004
$ws.Cells.Item(429, 51).Value = 0.78300000000000003
$ws.Cells.Item(429, 52).Value = 0
$ws.Cells.Item(429, 53).Value = 0.26100000000000001
$ws.Cells.Item(429, 54).Value = 0.78300000000000003
$ws.Cells.Item(429, 55).Value = 3.13
$ws.Cells.Item(429, 56).Value = 1.0429999999999999
$ws.Cells.Item(429, 57).Value = 0.64500000000000002
$ws.Cells.Item(429, 58).Value = 0.51600000000000001
$ws.Cells.Item(429, 59).Value = 1.032
$ws.Cells.Item(429, 60).Value = 1.806

$ws.Cells.Item(430, 1).Value = "EF611830"
$ws.Cells.Item(430, 2).Value = 0.74099999999999999
$ws.Cells.Item(430, 3).Value = 1.2589999999999999
$ws.Cells.Item(430, 4).Value = 0.68600000000000005
$ws.Cells.Item(430, 5).Value = 1.5429999999999999
$ws.Cells.Item(430, 6).Value = 0.51400000000000001
$ws.Cells.Item(430, 7).Value = 0.51400000000000001
$ws.Cells.Item(430, 8).Value = 0.85699999999999998
$ws.Cells.Item(430, 9).Value = 1.8859999999999999
$ws.Cells.Item(430, 10).Value = 1
$ws.Cells.Item(430, 11).Value = 1
$ws.Cells.Item(430, 12).Value = 1
$ws.Cells.Item(430, 13).Value = 0.88900000000000001
$ws.Cells.Item(430, 14).Value = 1.333
$ws.Cells.Item(430, 15).Value = 0.88900000000000001
$ws.Cells.Item(430, 16).Value = 0.88900000000000001
$ws.Cells.Item(430, 17).Value = 1.2
$ws.Cells.Item(430, 18).Value = 0.85699999999999998
$ws.Cells.Item(430, 19).Value = 1.714
$ws.Cells.Item(430, 20).Value = 0.68600000000000005
$ws.Cells.Item(430, 21).Value = 0.68600000000000005
$ws.Cells.Item(430, 22).Value = 0.85699999999999998
$ws.Cells.Item(430, 23).Value = 1.647
$ws.Cells.Item(430, 24).Value = 1.1759999999999999
$ws.Cells.Item(430, 25).Value = 1.1759999999999999
$ws.Cells.Item(430, 26).Value = 0
$ws.Cells.Item(430, 27).Value = 1.0669999999999999
$ws.Cells.Item(430, 28).Value = 1.0669999999999999
$ws.Cells.Item(430, 29).Value = 0.66700000000000004
$ws.Cells.Item(430, 30).Value = 1.2
$ws.Cells.Item(430, 31).Value = 1.4550000000000001
$ws.Cells.Item(430, 32).Value = 0.72699999999999998
$ws.Cells.Item(430, 33).Value = 1.6970000000000001
$ws.Cells.Item(430, 34).Value = 0.121
$ws.Cells.Item(430, 35).Value = 1.429
$ws.Cells.Item(430, 36).Value = 0.57099999999999995
$ws.Cells.Item(430, 37).Value = 0.92300000000000004
$ws.Cells.Item(430, 38).Value = 1.077
$ws.Cells.Item(430, 39).Value = 1.2
$ws.Cells.Item(430, 40).Value = 0.8
$ws.Cells.Item(430, 41).Value = 1
$ws.Cells.Item(430, 42).Value = 1
$ws.Cells.Item(430, 43).Value = 0.46200000000000002
$ws.Cells.Item(430, 44).Value = 1.538
$ws.Cells.Item(430, 45).Value = 1.0429999999999999
$ws.Cells.Item(430, 46).Value = 0.95699999999999996
$ws.Cells.Item(430, 47).Value = 0.56200000000000006
$ws.Cells.Item(430, 48).Value = 1.4379999999999999
$ws.Cells.Item(430, 49).Value = 1.333
$ws.Cells.Item(430, 50).Value = 0.66700000000000004
$ws.Cells.Item(430, 51).Value = 0.78300000000000003
$ws.Cells.Item(430, 52).Value = 0
$ws.Cells.Item(430, 53).Value = 0.26100000000000001
$ws.Cells.Item(430, 54).Value = 0.78300000000000003
$ws.Cells.Item(430, 55).Value = 3.13
$ws.Cells.Item(430, 56).Value = 1.0429999999999999
$ws.Cells.Item(430, 57).Value = 0.64500000000000002
$ws.Cells.Item(430, 58).Value = 0.51600000000000001
$ws.Cells.Item(430, 59).Value = 0.90300000000000002
$ws.Cells.Item(430, 60).Value = 1.9350000000000001

$ws.Cells.Item(431, 1).Value = "EF611851"
$ws.Cells.Item(431, 2).Value = 0.74099999999999999
$ws.Cells.Item(431, 3).Value = 1.2589999999999999
$ws.Cells.Item(431, 4).Value = 0.34300000000000003
$ws.Cells.Item(431, 5).Value = 1.714
$ws.Cells.Item(431, 6).Value = 0.68600000000000005
$ws.Cells.Item(431, 7).Value = 0.51400000000000001
$ws.Cells.Item(431, 8).Value = 1.0289999999999999
$ws.Cells.Item(431, 9).Value = 1.714
$ws.Cells.Item(431, 10).Value = 1.111
$ws.Cells.Item(431, 11).Value = 0.88900000000000001
$ws.Cells.Item(431, 12).Value = 1
$ws.Cells.Item(431, 13).Value = 1.1850000000000001
$ws.Cells.Item(431, 14).Value = 1.0369999999999999
$ws.Cells.Item(431, 15).Value = 0.88900000000000001
$ws.Cells.Item(431, 16).Value = 0.88900000000000001
$ws.Cells.Item(431, 17).Value = 1.371
$ws.Cells.Item(431, 18).Value = 0.68600000000000005
$ws.Cells.Item(431, 19).Value = 1.714
$ws.Cells.Item(431, 20).Value = 0.68600000000000005
$ws.Cells.Item(431, 21).Value = 0.85699999999999998
$ws.Cells.Item(431, 22).Value = 0.68600000000000005
$ws.Cells.Item(431, 23).Value = 1.647
$ws.Cells.Item(431, 24).Value = 1.1759999999999999
$ws.Cells.Item(431, 25).Value = 1.1759999999999999
$ws.Cells.Item(431, 26).Value = 0
$ws.Cells.Item(431, 27).Value = 1.6
$ws.Cells.Item(431, 28).Value = 0.66700000000000004
$ws.Cells.Item(431, 29).Value = 0.66700000000000004
$ws.Cells.Item(431, 30).Value = 1.0669999999999999
$ws.Cells.Item(431, 31).Value = 1.333
$ws.Cells.Item(431, 32).Value = 0.84799999999999998
$ws.Cells.Item(431, 33).Value = 1.6970000000000001
$ws.Cells.Item(431, 34).Value = 0.121
$ws.Cells.Item(431, 35).Value = 1.524
$ws.Cells.Item(431, 36).Value = 0.47599999999999998
$ws.Cells.Item(431, 37).Value = 0.92300000000000004
$ws.Cells.Item(431, 38).Value = 1.077
$ws.Cells.Item(431, 39).Value = 1.2
$ws.Cells.Item(431, 40).Value = 0.8
$ws.Cells.Item(431, 41).Value = 1.2
$ws.Cells.Item(431, 42).Value = 0.8
$ws.Cells.Item(431, 43).Value = 0.38500000000000001
$ws.Cells.Item(431, 44).Value = 1.615
$ws.Cells.Item(431, 45).Value = 1.0429999999999999
$ws.Cells.Item(431, 46).Value = 0.95699999999999996
$ws.Cells.Item(431, 47).Value = 0.68799999999999994
$ws.Cells.Item(431, 48).Value = 1.3120000000000001
$ws.Cells.Item(431, 49).Value = 1.333
$ws.Cells.Item(431, 50).Value = 0.66700000000000004
$ws.Cells.Item(431, 51).Value = 0.78300000000000003
$ws.Cells.Item(431, 52).Value = 0
$ws.Cells.Item(431, 53).Value = 0.26100000000000001
$ws.Cells.Item(431, 54).Value = 0.78300000000000003
$ws.Cells.Item(431, 55).Value = 3.391
$ws.Cells.Item(431, 56).Value = 0.78300000000000003
$ws.Cells.Item(431, 57).Value = 0.90300000000000002
$ws.Cells.Item(431, 58).Value = 0.38700000000000001
$ws.Cells.Item(431, 59).Value = 1.161
$ws.Cells.Item(431, 60).Value = 1.548

$ws.Cells.Item(432, 1).Value = "L20675"
$ws.Cells.Item(432, 2).Value = 0.74099999999999999
$ws.Cells.Item(432, 3).Value = 1.2589999999999999
$ws.Cells.Item(432, 4).Value = 0.68600000000000005
$ws.Cells.Item(432, 5).Value = 1.5429999999999999
$ws.Cells.Item(432, 6).Value = 0.51400000000000001
$ws.Cells.Item(432, 7).Value = 0.51400000000000001
$ws.Cells.Item(432, 8).Value = 0.68600000000000005
$ws.Cells.Item(432, 9).Value = 2.0569999999999999
$ws.Cells.Item(432, 10).Value = 1.222
$ws.Cells.Item(432, 11).Value = 0.77800000000000002
$ws.Cells.Item(432, 12).Value = 1
$ws.Cells.Item(432, 13).Value = 1.0369999999999999
$ws.Cells.Item(432, 14).Value = 1.1850000000000001
$ws.Cells.Item(432, 15).Value = 0.88900000000000001
$ws.Cells.Item(432, 16).Value = 0.88900000000000001
$ws.Cells.Item(432, 17).Value = 1.2350000000000001
$ws.Cells.Item(432, 18).Value = 1.0589999999999999
$ws.Cells.Item(432, 19).Value = 1.5880000000000001
$ws.Cells.Item(432, 20).Value = 0.52900000000000003
$ws.Cells.Item(432, 21).Value = 0.70599999999999996
$ws.Cells.Item(432, 22).Value = 0.88200000000000001
$ws.Cells.Item(432, 23).Value = 1.556
$ws.Cells.Item(432, 24).Value = 1.111
$ws.Cells.Item(432, 25).Value = 1.111
$ws.Cells.Item(432, 26).Value = 0.222
$ws.Cells.Item(432, 27).Value = 1.379
$ws.Cells.Item(432, 28).Value = 0.82799999999999996
$ws.Cells.Item(432, 29).Value = 0.82799999999999996
$ws.Cells.Item(432, 30).Value = 0.96599999999999997
$ws.Cells.Item(432, 31).Value = 1.294
$ws.Cells.Item(432, 32).Value = 0.82399999999999995
$ws.Cells.Item(432, 33).Value = 1.647
$ws.Cells.Item(432, 34).Value = 0.23499999999999999
$ws.Cells.Item(432, 35).Value = 1.333
$ws.Cells.Item(432, 36).Value = 0.66700000000000004
$ws.Cells.Item(432, 37).Value = 0.92300000000000004
$ws.Cells.Item(432, 38).Value = 1.077
$ws.Cells.Item(432, 39).Value = 1.2
$ws.Cells.Item(432, 40).Value = 0.8
$ws.Cells.Item(432, 41).Value = 0.85699999999999998
$ws.Cells.Item(432, 42).Value = 1.143
$ws.Cells.Item(432, 43).Value = 0.56000000000000005
$ws.Cells.Item(432, 44).Value = 1.44
$ws.Cells.Item(432, 45).Value = 1.0429999999999999
$ws.Cells.Item(432, 46).Value = 0.95699999999999996
$ws.Cells.Item(432, 47).Value = 0.68799999999999994
$ws.Cells.Item(432, 48).Value = 1.3120000000000001
$ws.Cells.Item(432, 49).Value = 1.333
$ws.Cells.Item(432, 50).Value = 0.66700000000000004
$ws.Cells.Item(432, 51).Value = 0.78300000000000003
$ws.Cells.Item(432, 52).Value = 0
$ws.Cells.Item(432, 53).Value = 0.26100000000000001
$ws.Cells.Item(432, 54).Value = 0.78300000000000003
$ws.Cells.Item(432, 55).Value = 3.13
$ws.Cells.Item(432, 56).Value = 1.0429999999999999
$ws.Cells.Item(432, 57).Value = 0.64500000000000002
$ws.Cells.Item(432, 58).Value = 0.38700000000000001
$ws.Cells.Item(432, 59).Value = 1.29
$ws.Cells.Item(432, 60).Value = 1.677

# Update the active selection to match the final view state after the edit
$ws.Range("F429").Select()
